$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2-23 from 45233 (2023-11-03) to 45243 (2023-11-13)
$ws.Range("C2:C23").Value = 45243
